# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# and sets the "Priority" column to "ht" for the rows that were (re)handed off.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows (on all three sheets) that correspond to the files included in this handoff.
$rows = @(7, 8, 9, 10, 11, 14)

foreach ($r in $rows) {
    # Overview sheet: "Latest HO Xliff Generate Date" column (G)
    $overview.Range("G$r").Value = "2016-08-13 14:24:54"

    # zh-cn sheet: "Latest Handoff Datetime" column (H) and "Priority" column (E)
    $zhcn.Range("H$r").Value = "2016-08-13 14:24:45"
    $zhcn.Range("E$r").Value = "ht"

    # de-de sheet: "Latest Handoff Datetime" column (H) and "Priority" column (E)
    $dede.Range("H$r").Value = "2016-08-13 14:24:54"
    $dede.Range("E$r").Value = "ht"
}
